$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A598").Value = "'2022-04-19"
$ws.Range("B598").Value = "overview"
$ws.Range("C598").Value = "K02000001"
$ws.Range("D598").Value = "United Kingdom"
$ws.Range("E598").Value = 21863944
$ws.Range("F598").Value = 116227
$ws.Range("G598").Value = 482
$ws.Range("H598").Value = 171878

$ws.Range("A599").Value = "'2022-04-20"
$ws.Range("B599").Value = "overview"
$ws.Range("C599").Value = "K02000001"
$ws.Range("D599").Value = "United Kingdom"
$ws.Range("E599").Value = 21890037
$ws.Range("F599").Value = 26147
$ws.Range("G599").Value = 508
$ws.Range("H599").Value = 172386

$ws.Range("A600").Value = "'2022-04-21"
$ws.Range("B600").Value = "overview"
$ws.Range("C600").Value = "K02000001"
$ws.Range("D600").Value = "United Kingdom"
$ws.Range("E600").Value = 21909509
$ws.Range("F600").Value = 19482
$ws.Range("G600").Value = 646
$ws.Range("H600").Value = 173032

$ws.Range("A601").Value = "'2022-04-22"
$ws.Range("B601").Value = "overview"
$ws.Range("C601").Value = "K02000001"
$ws.Range("D601").Value = "United Kingdom"
$ws.Range("E601").Value = 21933206
$ws.Range("F601").Value = 19795
$ws.Range("G601").Value = 284
$ws.Range("H601").Value = 173352

$ws.Range("A602").Value = "'2022-04-25"
$ws.Range("B602").Value = "overview"
$ws.Range("C602").Value = "K02000001"
$ws.Range("D602").Value = "United Kingdom"
$ws.Range("E602").Value = 21978198
$ws.Range("F602").Value = 45077
$ws.Range("G602").Value = 341
$ws.Range("H602").Value = 173693

$ws.Range("A603").Value = "'2022-04-26"
$ws.Range("B603").Value = "overview"
$ws.Range("C603").Value = "K02000001"
$ws.Range("D603").Value = "United Kingdom"
$ws.Range("E603").Value = 21994809
$ws.Range("F603").Value = 16579
$ws.Range("G603").Value = 451
$ws.Range("H603").Value = 174144

$ws.Range("A604").Value = "'2022-04-27"
$ws.Range("B604").Value = "overview"
$ws.Range("C604").Value = "K02000001"
$ws.Range("D604").Value = "United Kingdom"
$ws.Range("E604").Value = 22011920
$ws.Range("F604").Value = 17224
$ws.Range("G604").Value = 304
$ws.Range("H604").Value = 174448

$ws.Range("A605").Value = "'2022-04-28"
$ws.Range("B605").Value = "overview"
$ws.Range("C605").Value = "K02000001"
$ws.Range("D605").Value = "United Kingdom"
$ws.Range("E605").Value = 22025925
$ws.Range("F605").Value = 14030
$ws.Range("G605").Value = 248
$ws.Range("H605").Value = 174696

$ws.Range("A606").Value = "'2022-04-29"
$ws.Range("B606").Value = "overview"
$ws.Range("C606").Value = "K02000001"
$ws.Range("D606").Value = "United Kingdom"
$ws.Range("E606").Value = 22038340
$ws.Range("F606").Value = 12421
$ws.Range("G606").Value = 216
$ws.Range("H606").Value = 174912
